$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -3
$ws.Range("F5").Value = -1
$ws.Range("F8").Value = -1
